$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(375, 426.3756408691406),
    @(187.5, 148.458251953125),
    @(3375, 132.2644958496094),
    @(3703.125, 121.7497406005859),
    @(843.75, 105.3359527587891),
    @(1500, 97.34624481201172),
    @(937.5, 91.82764434814453),
    @(1218.75, 88.96601104736328),
    @(1406.25, 83.01222991943359),
    @(562.5, 78.421142578125)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
